# Auto-generated edit script: updates cached market-price / profit
# values (columns H-N) across multiple leve-profit sheets to match
# a refreshed market-data snapshot. Values were recomputed by the
# scheduled runner; this script replays the resulting cell writes.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 9138.852000000001  # H62: was 9531.208000000001
$ws.Cells.Item(62, 9).Value = 8075.4165  # I62: was 8767.223
$ws.Cells.Item(62, 11).Value = 8075.4165  # K62: was 8767.223
$ws.Cells.Item(62, 13).Value = -7451.4165  # M62: was -8143.223
$ws.Cells.Item(64, 8).Value = 7889.4  # H64: was 8022.8
$ws.Cells.Item(64, 9).Value = 6936.75  # I64: was 8582.666999999999
$ws.Cells.Item(64, 10).Value = 8524.5  # J64: was 7782.857
$ws.Cells.Item(64, 11).Value = 6936.75  # K64: was 8582.666999999999
$ws.Cells.Item(64, 12).Value = 8524.5  # L64: was 7782.857
$ws.Cells.Item(64, 13).Value = -6688.75  # M64: was -8334.666999999999
$ws.Cells.Item(64, 14).Value = -9020.5  # N64: was -8278.857
$ws.Cells.Item(65, 8).Value = 9138.852000000001  # H65: was 9531.208000000001
$ws.Cells.Item(65, 9).Value = 8075.4165  # I65: was 8767.223
$ws.Cells.Item(65, 11).Value = 40377.0825  # K65: was 43836.115
$ws.Cells.Item(65, 13).Value = -37257.0825  # M65: was -40716.115
$ws.Cells.Item(67, 8).Value = 7889.4  # H67: was 8022.8
$ws.Cells.Item(67, 9).Value = 6936.75  # I67: was 8582.666999999999
$ws.Cells.Item(67, 10).Value = 8524.5  # J67: was 7782.857
$ws.Cells.Item(67, 11).Value = 6936.75  # K67: was 8582.666999999999
$ws.Cells.Item(67, 12).Value = 8524.5  # L67: was 7782.857
$ws.Cells.Item(67, 13).Value = -6078.75  # M67: was -7724.666999999999
$ws.Cells.Item(67, 14).Value = -10240.5  # N67: was -9498.857
$ws.Cells.Item(70, 8).Value = 3843.6667  # H70: was 4316.3335
$ws.Cells.Item(70, 9).Value = 2599.6667  # I70: was 2999
$ws.Cells.Item(70, 10).Value = 4465.6665  # J70: was 4579.8
$ws.Cells.Item(70, 11).Value = 7799.000100000001  # K70: was 8997
$ws.Cells.Item(70, 12).Value = 13396.9995  # L70: was 13739.4
$ws.Cells.Item(70, 13).Value = -7529.000100000001  # M70: was -8727
$ws.Cells.Item(70, 14).Value = -13936.9995  # N70: was -14279.4
$ws.Cells.Item(73, 8).Value = 3843.6667  # H73: was 4316.3335
$ws.Cells.Item(73, 9).Value = 2599.6667  # I73: was 2999
$ws.Cells.Item(73, 10).Value = 4465.6665  # J73: was 4579.8
$ws.Cells.Item(73, 11).Value = 7799.000100000001  # K73: was 8997
$ws.Cells.Item(73, 12).Value = 13396.9995  # L73: was 13739.4
$ws.Cells.Item(73, 13).Value = -6863.000100000001  # M73: was -8061
$ws.Cells.Item(73, 14).Value = -15268.9995  # N73: was -15611.4
$ws.Cells.Item(86, 8).Value = 4206.6  # H86: was 3799.9443
$ws.Cells.Item(86, 9).Value = 3563.5454  # I86: was 3516.5833
$ws.Cells.Item(86, 10).Value = 5975  # J86: was 4366.6665
$ws.Cells.Item(86, 11).Value = 3563.5454  # K86: was 3516.5833
$ws.Cells.Item(86, 12).Value = 5975  # L86: was 4366.6665
$ws.Cells.Item(86, 13).Value = -2440.5454  # M86: was -2393.5833
$ws.Cells.Item(86, 14).Value = -8221  # N86: was -6612.6665
$ws.Cells.Item(89, 8).Value = 4206.6  # H89: was 3799.9443
$ws.Cells.Item(89, 9).Value = 3563.5454  # I89: was 3516.5833
$ws.Cells.Item(89, 10).Value = 5975  # J89: was 4366.6665
$ws.Cells.Item(89, 11).Value = 17817.727  # K89: was 17582.9165
$ws.Cells.Item(89, 12).Value = 29875  # L89: was 21833.3325
$ws.Cells.Item(89, 13).Value = -12201.727  # M89: was -11966.9165
$ws.Cells.Item(89, 14).Value = -41107  # N89: was -33065.3325
$ws.Cells.Item(95, 8).Value = 0  # H95: was 60000
$ws.Cells.Item(95, 10).Value = 0  # J95: was 60000
$ws.Cells.Item(95, 12).Value = 0  # L95: was 60000
$ws.Cells.Item(95, 14).ClearContents()  # N95: was -65492
$ws.Cells.Item(125, 8).Value = 1437.1666  # H125: was 1473
$ws.Cells.Item(125, 9).Value = 1233.5  # I125: was 1234
$ws.Cells.Item(125, 10).Value = 1539  # J125: was 1592.5
$ws.Cells.Item(125, 11).Value = 11101.5  # K125: was 11106
$ws.Cells.Item(125, 12).Value = 13851  # L125: was 14332.5
$ws.Cells.Item(125, 13).Value = -8641.5  # M125: was -8646
$ws.Cells.Item(125, 14).Value = -18771  # N125: was -19252.5
$ws.Cells.Item(132, 8).Value = 8179833.5  # H132: was 9315904
$ws.Cells.Item(132, 9).Value = 9834376  # I132: was 10786106
$ws.Cells.Item(132, 10).Value = 143485  # J132: was 200649.8
$ws.Cells.Item(132, 11).Value = 29503128  # K132: was 32358318
$ws.Cells.Item(132, 12).Value = 430455  # L132: was 601949.3999999999
$ws.Cells.Item(132, 13).Value = -29500598  # M132: was -32355788
$ws.Cells.Item(132, 14).Value = -435515  # N132: was -607009.3999999999
$ws.Cells.Item(137, 8).Value = 2330  # H137: was 2422.125
$ws.Cells.Item(137, 9).Value = 2648.9  # I137: was 2799
$ws.Cells.Item(137, 10).Value = 1931.375  # J137: was 1937.5714
$ws.Cells.Item(137, 11).Value = 7946.700000000001  # K137: was 8397
$ws.Cells.Item(137, 12).Value = 5794.125  # L137: was 5812.7142
$ws.Cells.Item(137, 13).Value = -5396.700000000001  # M137: was -5847
$ws.Cells.Item(137, 14).Value = -10894.125  # N137: was -10912.7142
$ws.Cells.Item(138, 8).Value = 3632.7568  # H138: was 3595.8027
$ws.Cells.Item(138, 9).Value = 1503.8667  # I138: was 1471.5483
$ws.Cells.Item(138, 10).Value = 5084.273  # J138: was 5059.1777
$ws.Cells.Item(138, 11).Value = 4511.6001  # K138: was 4414.644899999999
$ws.Cells.Item(138, 12).Value = 15252.819  # L138: was 15177.5331
$ws.Cells.Item(138, 13).Value = 628.3999000000003  # M138: was 725.3551000000007
$ws.Cells.Item(138, 14).Value = -25532.819  # N138: was -25457.5331

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4072.9756  # H32: was 4276.795
$ws.Cells.Item(32, 9).Value = 3202.7878  # I32: was 3403.0645
$ws.Cells.Item(32, 11).Value = 3202.7878  # K32: was 3403.0645
$ws.Cells.Item(32, 13).Value = -2915.7878  # M32: was -3116.0645
$ws.Cells.Item(45, 8).Value = 12308.454  # H45: was 12390.091
$ws.Cells.Item(45, 9).Value = 17285.572  # I45: was 17413.857
$ws.Cells.Item(45, 11).Value = 17285.572  # K45: was 17413.857
$ws.Cells.Item(45, 13).Value = -16908.572  # M45: was -17036.857
$ws.Cells.Item(74, 8).Value = 1736.3334  # H74: was 1874.5
$ws.Cells.Item(74, 9).Value = 1387.5714  # I74: was 1524.4
$ws.Cells.Item(74, 11).Value = 1387.5714  # K74: was 1524.4
$ws.Cells.Item(74, 13).Value = -513.5714  # M74: was -650.4000000000001
$ws.Cells.Item(77, 8).Value = 1736.3334  # H77: was 1874.5
$ws.Cells.Item(77, 9).Value = 1387.5714  # I77: was 1524.4
$ws.Cells.Item(77, 11).Value = 6937.857  # K77: was 7622
$ws.Cells.Item(77, 13).Value = -2569.857  # M77: was -3254
$ws.Cells.Item(122, 8).Value = 113084.445  # H122: was 126468.375
$ws.Cells.Item(122, 10).Value = 6013.5  # J122: was 6014
$ws.Cells.Item(122, 12).Value = 18040.5  # L122: was 18042
$ws.Cells.Item(122, 14).Value = -22940.5  # N122: was -22942

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 3183.158  # H107: was 3249
$ws.Cells.Item(107, 10).Value = 1999.4  # J107: was 1999.75
$ws.Cells.Item(107, 12).Value = 1999.4  # L107: was 1999.75
$ws.Cells.Item(107, 14).Value = -5839.4  # N107: was -5839.75

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 60390.832  # H31: was 57319
$ws.Cells.Item(31, 9).Value = 72894.36  # I31: was 84962.664
$ws.Cells.Item(31, 10).Value = 16628.5  # J31: was 9929.857
$ws.Cells.Item(31, 11).Value = 72894.36  # K31: was 84962.664
$ws.Cells.Item(31, 12).Value = 16628.5  # L31: was 9929.857
$ws.Cells.Item(31, 13).Value = -72599.36  # M31: was -84667.664
$ws.Cells.Item(31, 14).Value = -17218.5  # N31: was -10519.857
$ws.Cells.Item(34, 8).Value = 60390.832  # H34: was 57319
$ws.Cells.Item(34, 9).Value = 72894.36  # I34: was 84962.664
$ws.Cells.Item(34, 10).Value = 16628.5  # J34: was 9929.857
$ws.Cells.Item(34, 11).Value = 72894.36  # K34: was 84962.664
$ws.Cells.Item(34, 12).Value = 16628.5  # L34: was 9929.857
$ws.Cells.Item(34, 13).Value = -72692.36  # M34: was -84760.664
$ws.Cells.Item(34, 14).Value = -17032.5  # N34: was -10333.857
$ws.Cells.Item(99, 8).Value = 2788.3572  # H99: was 2959.1667
$ws.Cells.Item(99, 9).Value = 2501  # I99: was 2610.6667
$ws.Cells.Item(99, 10).Value = 3506.75  # J99: was 4004.6667
$ws.Cells.Item(99, 11).Value = 2501  # K99: was 2610.6667
$ws.Cells.Item(99, 12).Value = 3506.75  # L99: was 4004.6667
$ws.Cells.Item(99, 13).Value = -1003  # M99: was -1112.6667
$ws.Cells.Item(99, 14).Value = -6502.75  # N99: was -7000.6667
$ws.Cells.Item(105, 9).Value = 1564.75  # I105: was 1496.4445
$ws.Cells.Item(105, 10).Value = 1480.5  # J105: was 2011
$ws.Cells.Item(105, 11).Value = 1564.75  # K105: was 1496.4445
$ws.Cells.Item(105, 12).Value = 1480.5  # L105: was 2011
$ws.Cells.Item(105, 13).Value = 182.25  # M105: was 250.5554999999999
$ws.Cells.Item(105, 14).Value = -4974.5  # N105: was -5505
$ws.Cells.Item(107, 8).Value = 1324.1428  # H107: was 1400.2307
$ws.Cells.Item(107, 9).Value = 1359.7  # I107: was 1403.7
$ws.Cells.Item(107, 10).Value = 1235.25  # J107: was 1388.6666
$ws.Cells.Item(107, 11).Value = 1359.7  # K107: was 1403.7
$ws.Cells.Item(107, 12).Value = 1235.25  # L107: was 1388.6666
$ws.Cells.Item(107, 13).Value = 560.3  # M107: was 516.3
$ws.Cells.Item(107, 14).Value = -5075.25  # N107: was -5228.6666
$ws.Cells.Item(126, 8).Value = 2788.3572  # H126: was 2959.1667
$ws.Cells.Item(126, 9).Value = 2501  # I126: was 2610.6667
$ws.Cells.Item(126, 10).Value = 3506.75  # J126: was 4004.6667
$ws.Cells.Item(126, 11).Value = 7503  # K126: was 7832.000100000001
$ws.Cells.Item(126, 12).Value = 10520.25  # L126: was 12014.0001
$ws.Cells.Item(126, 13).Value = -5033  # M126: was -5362.000100000001
$ws.Cells.Item(126, 14).Value = -15460.25  # N126: was -16954.0001
$ws.Cells.Item(132, 8).Value = 3299.0393  # H132: was 3391.9792
$ws.Cells.Item(132, 9).Value = 3254.9268  # I132: was 3333.7437
$ws.Cells.Item(132, 10).Value = 3479.9  # J132: was 3644.3333
$ws.Cells.Item(132, 11).Value = 9764.7804  # K132: was 10001.2311
$ws.Cells.Item(132, 12).Value = 10439.7  # L132: was 10932.9999
$ws.Cells.Item(132, 13).Value = -7234.7804  # M132: was -7471.231100000001
$ws.Cells.Item(132, 14).Value = -15499.7  # N132: was -15992.9999
$ws.Cells.Item(134, 8).Value = 7034.2095  # H134: was 7036.613
$ws.Cells.Item(134, 9).Value = 5085.8364  # I134: was 5088.5454
$ws.Cells.Item(134, 11).Value = 15257.5092  # K134: was 15265.6362
$ws.Cells.Item(134, 13).Value = -12722.5092  # M134: was -12730.6362

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 31311320  # H4: was 24083554
$ws.Cells.Item(4, 9).Value = 36799816  # I4: was 25533588
$ws.Cells.Item(4, 10).Value = 209833.33  # J4: was 399666.66
$ws.Cells.Item(4, 11).Value = 110399448  # K4: was 76600764
$ws.Cells.Item(4, 12).Value = 629499.99  # L4: was 1198999.98
$ws.Cells.Item(4, 13).Value = -110399336  # M4: was -76600652
$ws.Cells.Item(4, 14).Value = -629723.99  # N4: was -1199223.98
$ws.Cells.Item(80, 8).Value = 3824.75  # H80: was 3859.8
$ws.Cells.Item(80, 10).Value = 3850  # J80: was 3900
$ws.Cells.Item(80, 12).Value = 11550  # L80: was 11700
$ws.Cells.Item(80, 14).Value = -13422  # N80: was -13572
$ws.Cells.Item(83, 8).Value = 3824.75  # H83: was 3859.8
$ws.Cells.Item(83, 10).Value = 3850  # J83: was 3900
$ws.Cells.Item(83, 12).Value = 34650  # L83: was 35100
$ws.Cells.Item(83, 14).Value = -44010  # N83: was -44460
$ws.Cells.Item(129, 8).Value = 840.1667  # H129: was 920.1429000000001
$ws.Cells.Item(129, 10).Value = 1500  # J129: was 1450
$ws.Cells.Item(129, 12).Value = 4500  # L129: was 4350
$ws.Cells.Item(129, 14).Value = -14500  # N129: was -14350
$ws.Cells.Item(131, 8).Value = 159637.36  # H131: was 75489.47
$ws.Cells.Item(131, 9).Value = 500572.5  # I131: was 667053.3
$ws.Cells.Item(131, 10).Value = 23263.3  # J131: was 9760.147999999999
$ws.Cells.Item(131, 11).Value = 1501717.5  # K131: was 2001159.9
$ws.Cells.Item(131, 12).Value = 69789.89999999999  # L131: was 29280.444
$ws.Cells.Item(131, 13).Value = -1496677.5  # M131: was -1996119.9
$ws.Cells.Item(131, 14).Value = -79869.89999999999  # N131: was -39360.444
$ws.Cells.Item(137, 8).Value = 3136.3333  # H137: was 3228.1428
$ws.Cells.Item(137, 9).Value = 2921.1667  # I137: was 2939.4
$ws.Cells.Item(137, 10).Value = 3566.6667  # J137: was 3950
$ws.Cells.Item(137, 11).Value = 8763.500100000001  # K137: was 8818.200000000001
$ws.Cells.Item(137, 12).Value = 10700.0001  # L137: was 11850
$ws.Cells.Item(137, 13).Value = -3663.500100000001  # M137: was -3718.200000000001
$ws.Cells.Item(137, 14).Value = -20900.0001  # N137: was -22050
$ws.Cells.Item(140, 8).Value = 4004.5715  # H140: was 4528.25
$ws.Cells.Item(140, 9).Value = 3906.4  # I140: was 4528.25
$ws.Cells.Item(140, 10).Value = 4250  # J140: was 0
$ws.Cells.Item(140, 11).Value = 11719.2  # K140: was 13584.75
$ws.Cells.Item(140, 12).Value = 12750  # L140: was 0
$ws.Cells.Item(140, 13).Value = -6539.200000000001  # M140: was -8404.75
$ws.Cells.Item(140, 14).Value = -23110  # N140: was None

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 2199.48  # H107: was 2388.7827
$ws.Cells.Item(107, 9).Value = 1779.0667  # I107: was 1924.2142
$ws.Cells.Item(107, 10).Value = 2830.1  # J107: was 3111.4443
$ws.Cells.Item(107, 11).Value = 1779.0667  # K107: was 1924.2142
$ws.Cells.Item(107, 12).Value = 2830.1  # L107: was 3111.4443
$ws.Cells.Item(107, 13).Value = 140.9332999999999  # M107: was -4.214199999999892
$ws.Cells.Item(107, 14).Value = -6670.1  # N107: was -6951.4443
$ws.Cells.Item(122, 8).Value = 3714.1365  # H122: was 3867.9
$ws.Cells.Item(122, 9).Value = 3397.5  # I122: was 3477.8462
$ws.Cells.Item(122, 10).Value = 4268.25  # J122: was 4592.2856
$ws.Cells.Item(122, 11).Value = 10192.5  # K122: was 10433.5386
$ws.Cells.Item(122, 12).Value = 12804.75  # L122: was 13776.8568
$ws.Cells.Item(122, 13).Value = -7742.5  # M122: was -7983.5386
$ws.Cells.Item(122, 14).Value = -17704.75  # N122: was -18676.8568
$ws.Cells.Item(123, 8).Value = 21258.908  # H123: was 23994
$ws.Cells.Item(123, 10).Value = 21884.8  # J123: was 25118.25
$ws.Cells.Item(123, 12).Value = 21884.8  # L123: was 25118.25
$ws.Cells.Item(123, 14).Value = -26784.8  # N123: was -30018.25

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 7807.4116  # H7: was 6303.7827
$ws.Cells.Item(7, 9).Value = 7732.875  # I7: was 7035.6665
$ws.Cells.Item(7, 10).Value = 9000  # J7: was 3669
$ws.Cells.Item(7, 11).Value = 7732.875  # K7: was 7035.6665
$ws.Cells.Item(7, 12).Value = 9000  # L7: was 3669
$ws.Cells.Item(7, 13).Value = -7620.875  # M7: was -6923.6665
$ws.Cells.Item(7, 14).Value = -9224  # N7: was -3893
$ws.Cells.Item(40, 8).Value = 4821.1177  # H40: was 5186.1875
$ws.Cells.Item(40, 9).Value = 4073  # I40: was 4633.5454
$ws.Cells.Item(40, 10).Value = 7252.5  # J40: was 6402
$ws.Cells.Item(40, 11).Value = 4073  # K40: was 4633.5454
$ws.Cells.Item(40, 12).Value = 7252.5  # L40: was 6402
$ws.Cells.Item(40, 13).Value = -3937  # M40: was -4497.5454
$ws.Cells.Item(40, 14).Value = -7524.5  # N40: was -6674
$ws.Cells.Item(93, 8).Value = 4099.375  # H93: was 4470.857
$ws.Cells.Item(93, 9).Value = 3270.7144  # I93: was 3566
$ws.Cells.Item(93, 11).Value = 3270.7144  # K93: was 3566
$ws.Cells.Item(93, 13).Value = -2022.7144  # M93: was -2318
$ws.Cells.Item(126, 8).Value = 7807.4116  # H126: was 6303.7827
$ws.Cells.Item(126, 9).Value = 7732.875  # I126: was 7035.6665
$ws.Cells.Item(126, 10).Value = 9000  # J126: was 3669
$ws.Cells.Item(126, 11).Value = 23198.625  # K126: was 21106.9995
$ws.Cells.Item(126, 12).Value = 27000  # L126: was 11007
$ws.Cells.Item(126, 13).Value = -20728.625  # M126: was -18636.9995
$ws.Cells.Item(126, 14).Value = -31940  # N126: was -15947
$ws.Cells.Item(133, 8).Value = 48991.668  # H133: was 60000
$ws.Cells.Item(133, 10).Value = 48991.668  # J133: was 60000
$ws.Cells.Item(133, 12).Value = 48991.668  # L133: was 60000
$ws.Cells.Item(133, 14).Value = -54051.668  # N133: was -65060

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 1159.4517  # H100: was 1194.931
$ws.Cells.Item(100, 9).Value = 878.9375  # I100: was 894.2
$ws.Cells.Item(100, 10).Value = 1458.6666  # J100: was 1517.1428
$ws.Cells.Item(100, 11).Value = 1757.875  # K100: was 1788.4
$ws.Cells.Item(100, 12).Value = 2917.3332  # L100: was 3034.2856
$ws.Cells.Item(100, 13).Value = -1216.875  # M100: was -1247.4
$ws.Cells.Item(100, 14).Value = -3999.3332  # N100: was -4116.2856
$ws.Cells.Item(113, 8).Value = 530.7778  # H113: was 704.5
$ws.Cells.Item(113, 9).Value = 334.25  # I113: was 493.5
$ws.Cells.Item(113, 10).Value = 688  # J113: was 810
$ws.Cells.Item(113, 11).Value = 1002.75  # K113: was 1480.5
$ws.Cells.Item(113, 12).Value = 2064  # L113: was 2430
$ws.Cells.Item(113, 13).Value = 1167.25  # M113: was 689.5
$ws.Cells.Item(113, 14).Value = -6404  # N113: was -6770
$ws.Cells.Item(126, 8).Value = 4062.3333  # H126: was 5053.905
$ws.Cells.Item(126, 9).Value = 2609.75  # I126: was 2693.3635
$ws.Cells.Item(126, 10).Value = 6967.5  # J126: was 7650.5
$ws.Cells.Item(126, 11).Value = 7829.25  # K126: was 8080.0905
$ws.Cells.Item(126, 12).Value = 20902.5  # L126: was 22951.5
$ws.Cells.Item(126, 13).Value = -5359.25  # M126: was -5610.0905
$ws.Cells.Item(126, 14).Value = -25842.5  # N126: was -27891.5
$ws.Cells.Item(132, 8).Value = 6795.273  # H132: was 7674.75
$ws.Cells.Item(132, 9).Value = 6774.8  # I132: was 7771.143
$ws.Cells.Item(132, 11).Value = 20324.4  # K132: was 23313.429
$ws.Cells.Item(132, 13).Value = -17794.4  # M132: was -20783.429
$ws.Cells.Item(136, 8).Value = 2657.2856  # H136: was 2711.6296
$ws.Cells.Item(136, 9).Value = 2620.4092  # I136: was 2688.524
$ws.Cells.Item(136, 11).Value = 7861.2276  # K136: was 8065.572
$ws.Cells.Item(136, 13).Value = -5311.2276  # M136: was -5515.572

